$d = $word.ActiveDocument

# --- 1 & 2) "Going through the <spellcheck>riipen</spellcheck> project list ..."
#     (two identical occurrences in the document) — drop the spell-check
#     markup around "riipen" by collapsing the three runs into one.
$old1 = "Going through the riipen project list to select feasible project matching to our resource capability "
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $old1, 2)

# --- 3 & 4) " Discussion with the instructor on how to approach the <spellcheck>riipen</spellcheck> project."
#     (two identical occurrences) — same spell-check cleanup.
$old2 = " Discussion with the instructor on how to approach the riipen project."
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2)

# --- 5) "October 28, 2025 " — drop the grammar-check markup around "2025".
$old3 = "October 28, 2025 "
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $old3, 2)

# --- 6) "Updated Riipen dashboard to the latest milestone." — drop the
#     spell-check markup around "Riipen".
$old4 = "Updated Riipen dashboard to the latest milestone."
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $old4, 2)

# --- 7) "... as a survey, and then communicated ..." — drop the
#     grammar-check markup around "survey, and".
$old5 = "Created the Microsoft form to collect the user response as a survey, and then communicated to the other member about the progress, rectified changes and then "
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $old5, 2)

# --- 8) "Generate storyboards with AI (art form) for the personas (loosely
#     based) we wrote." -> mention the persona names "Ashley" and "Eleanor"
#     and add a trailing space.
$dash = [char]0x2013
$old6 = "Generate storyboards with AI (art form) for the personas (loosely based) we wrote."
$new6 = "Generate storyboards with AI (art form " + $dash + " Ashley and Eleanor) for the personas (loosely based) we wrote. "
$d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2)
